$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.3848

$ws.Range("A4").Value = -21.54830000000002
$ws.Range("B4").Value = 5.077899999999999
$ws.Range("E4").Value = 12.82530000000002

$ws.Range("B5").Value = 5.180999999999998

$ws.Range("A6").Value = -21.2947

$ws.Range("A7").Value = -21.4684

$ws.Range("B8").Value = 5.1172

$ws.Range("E9").Value = 13.3101

$ws.Range("E11").Value = 13.5876

$ws.Range("E14").Value = 12.6787

$ws.Range("A16").Value = -21.47400000000002
$ws.Range("B16").Value = 5.230799999999999

$ws.Range("E18").Value = 13.08220000000001

$ws.Range("A20").Value = -22.28630000000003

$ws.Range("B22").Value = 5.174400000000005

$ws.Range("E25").Value = 12.75560000000001
